# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to be stored as text so values such as
# "0.9989" or "0.000007996" keep their exact original formatting
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "26.696.16"
$ws.Range("E2").Value = "  +0.89%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.856.41"
$ws.Range("E3").Value = "  +0.60%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "0.9989"
$ws.Range("E4").Value = "  -0.37%  "

# Row 5 - BNB
$ws.Range("D5").Value = "265.63"
$ws.Range("E5").Value = "  +2.48%  "

# Row 6 - USDC
$ws.Range("D6").Value = "0.9995"
$ws.Range("E6").Value = "  -0.17%  "

# Row 7 - XRP
$ws.Range("D7").Value = "0.5242"
$ws.Range("E7").Value = "  +0.40%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "0.3295"
$ws.Range("E8").Value = "  +0.69%  "

# Row 9 - Dogecoin
$ws.Range("D9").Value = "0.06809"
$ws.Range("E9").Value = "  +1.04%  "

# Row 10 - Solana
$ws.Range("D10").Value = "18.91"
$ws.Range("E10").Value = "  -2.54%  "

# Row 11 - Polygon
$ws.Range("D11").Value = "0.7787"
$ws.Range("E11").Value = "  +0.58%  "

# Row 12 - TRON
$ws.Range("D12").Value = "0.07740"
$ws.Range("E12").Value = "  +0.92%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.894.14"
$ws.Range("E13").Value = "  +2.77%  "

# Row 14 - Litecoin
$ws.Range("D14").Value = "88.78"
$ws.Range("E14").Value = "  +0.32%  "

# Row 15 - Polkadot
$ws.Range("D15").Value = "5.043"
$ws.Range("E15").Value = "  -0.18%  "

# Row 16 - BinanceUSD
$ws.Range("D16").Value = "0.9983"
$ws.Range("E16").Value = "  -0.39%  "

# Row 17 - Avalanche
$ws.Range("D17").Value = "14.07"
$ws.Range("E17").Value = "  -0.50%  "

# Row 18 - ShibaInu
$ws.Range("D18").Value = "0.000007996"
$ws.Range("E18").Value = "  +1.12%  "

# Row 19 - Dai
$ws.Range("D19").Value = "0.9990"
$ws.Range("E19").Value = "  -0.29%  "

# Row 20 - WrappedBTC
$ws.Range("D20").Value = "26.695.46"
$ws.Range("E20").Value = "  +0.73%  "

# Row 21 - WrappedliquidstakedEther2.0
$ws.Range("D21").Value = "2.086.28"
$ws.Range("E21").Value = "  -0.57%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "4.641"
$ws.Range("E22").Value = "  +0.74%  "

# Row 23 - Cosmos
$ws.Range("D23").Value = "9.562"
$ws.Range("E23").Value = "  -0.47%  "

# Row 24 - Chainlink
$ws.Range("D24").Value = "6.019"
$ws.Range("E24").Value = "  +0.35%  "

# Row 25 - Monero
$ws.Range("D25").Value = "144.24"
$ws.Range("E25").Value = "  -0.64%  "

# Row 26 - LidoDAOToken
$ws.Range("D26").Value = "2.214"
$ws.Range("E26").Value = "  -5.04%  "

# Row 27 - Toncoin
$ws.Range("D27").Value = "1.681"
$ws.Range("E27").Value = "  +2.42%  "

# Row 28 - EthereumClassic
$ws.Range("D28").Value = "17.03"
$ws.Range("E28").Value = "  -0.06%  "

# Row 29 - BitcoinCash
$ws.Range("D29").Value = "112.38"
$ws.Range("E29").Value = "  +0.92%  "

# Row 30 - InternetComputer(DFINITY)
$ws.Range("D30").Value = "4.209"
$ws.Range("E30").Value = "  -0.58%  "

# Row 31 - Filecoin
$ws.Range("D31").Value = "4.165"
$ws.Range("E31").Value = "  -0.71%  "

# Row 32 - Stellar
$ws.Range("D32").Value = "0.08765"
$ws.Range("E32").Value = "  -0.01%  "

# Row 33 - Hedera
$ws.Range("D33").Value = "0.04842"
$ws.Range("E33").Value = "  -0.14%  "

# Row 34 - ARBITRUM
$ws.Range("D34").Value = "1.142"
$ws.Range("E34").Value = "  +0.15%  "

# Rows 35 & 36 - HuobiToken and ImmutableX swap positions
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "0.7172"
$ws.Range("E35").Value = "  +1.42%  "

$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Value = "2.868"
$ws.Range("E36").Value = "  -0.06%  "

# Row 37 - MXToken
$ws.Range("D37").Value = "3.119"
$ws.Range("E37").Value = "  +0.27%  "

# Row 38 - VeChain
$ws.Range("D38").Value = "0.01790"
$ws.Range("E38").Value = "  -1.27%  "

# Row 39 - RenderToken
$ws.Range("D39").Value = "2.213"
$ws.Range("E39").Value = "  -0.76%  "

# Row 40 - TheSandbox
$ws.Range("D40").Value = "0.4900"
$ws.Range("E40").Value = "  -0.93%  "

# Row 41 - Quant
$ws.Range("D41").Value = "112.96"
$ws.Range("E41").Value = "  +0.24%  "

# Row 42 - TrustWalletToken
$ws.Range("D42").Value = "0.9025"
$ws.Range("E42").Value = "  +0.11%  "

# Row 43 - FraxShare
$ws.Range("D43").Value = "6.092"
$ws.Range("E43").Value = "  +0.13%  "

# Row 44 - PaxDollar
$ws.Range("D44").Value = "0.9991"
$ws.Range("E44").Value = "  -0.29%  "

# Row 45 - Aptos
$ws.Range("D45").Value = "7.747"
$ws.Range("E45").Value = "  -0.73%  "

# Row 46 - Decentraland
$ws.Range("D46").Value = "0.4210"
$ws.Range("E46").Value = "  -1.80%  "

# Rows 47 & 48 - EnergySwap and Cronos swap positions
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "0.05926"
$ws.Range("E47").Value = "  +0.10%  "

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "9.098"
$ws.Range("E48").Value = "  -1.39%  "

# Row 49 - Algorand (only volume changes)
$ws.Range("E49").Value = "  -3.84%  "

# Row 50 - Elrond
$ws.Range("D50").Value = "35.13"
$ws.Range("E50").Value = "  -0.49%  "

# Row 51 - Aave
$ws.Range("D51").Value = "60.32"
$ws.Range("E51").Value = "  +1.11%  "

# Restore the default "Normal" style on the Price column so the
# underlying cell styling matches the original workbook (no stray
# explicit style index left behind from the text-format coercion above).
$ws.Range("D2:D51").Style = "Normal"
